# feat: support duplicate headers with non-unique keys
# Adds a new "Rebates-Purchases" worksheet (with monthly Rebates/Purchases
# columns per person, duplicate "Rebates"/"Purchases" headers) and tweaks
# the Departments sheet's selected cell.

$wb = $excel.ActiveWorkbook

# --- add the new worksheet at the end of the tab strip --------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Rebates-Purchases"

# --- row 1: paired month-end dates across C1:Z1 (one date feeds the
#     "Rebates" column and the next feeds "Purchases") -------------------
# Serial date values for the 23rd of each month, Jan-Dec 2023.
$dateSerials = @(44949,44980,45008,45039,45069,45100,45130,45161,45192,45222,45253,45283)
$col = 3
foreach ($serial in $dateSerials) {
  $ws.Cells.Item(1, $col).Value = $serial
  $ws.Cells.Item(1, $col).NumberFormat = "d-mmm"
  $ws.Cells.Item(1, $col + 1).Value = $serial
  $ws.Cells.Item(1, $col + 1).NumberFormat = "d-mmm"
  $col = $col + 2
}

# --- row 2: header labels, with "Rebates"/"Purchases" repeating for each
#     month (duplicate, non-unique keys) -----------------------------------
$ws.Cells.Item(2, 1).Value = "Name*"
$ws.Cells.Item(2, 2).Value = "Group*"
$col = 3
for ($i = 0; $i -lt 12; $i++) {
  $ws.Cells.Item(2, $col).Value = "Rebates"
  $ws.Cells.Item(2, $col + 1).Value = "Purchases"
  $col = $col + 2
}

# --- rows 3-7: one row per person, Rebates/Purchases growing by row and
#     by month ---------------------------------------------------------
$people = @(
  @{ Name = "John Doe";       Group = "Group A" },
  @{ Name = "Jane Smith";     Group = "Group B" },
  @{ Name = "David Johnson";  Group = "Group C" },
  @{ Name = "Lisa Adams";     Group = "Group D" },
  @{ Name = "Mary Johnson";   Group = "Group E" }
)

$row = 3
for ($r = 0; $r -lt $people.Count; $r++) {
  $ws.Cells.Item($row, 1).Value = $people[$r].Name
  $ws.Cells.Item($row, 2).Value = $people[$r].Group

  $col = 3
  for ($p = 0; $p -lt 12; $p++) {
    $rebate = 100 + $r * 100 + $p * 10
    $purchase = 1000 + $r * 1000 + $p * 100
    $ws.Cells.Item($row, $col).Value = $rebate
    $ws.Cells.Item($row, $col + 1).Value = $purchase
    $col = $col + 2
  }
  $row = $row + 1
}

# --- column A width (matches source workbook's best-fit width) -----------
$ws.Columns.Item(1).ColumnWidth = 12.1640625

# --- this sheet's own last-known selection -------------------------------
$ws.Range("B2").Select() | Out-Null

# --- restore Departments as the active tab, with its new selection -------
$dept = $wb.Worksheets.Item("Departments")
$dept.Activate()
$dept.Range("D5").Select() | Out-Null

# --- workbook window Y position tweak (best-effort; some hosts don't
#     expose this attribute through the object model) ---------------------
$wb.Windows.Item(1).Top = 1900
